$d = $word.ActiveDocument
$d.Content.Find.Execute("<id>p016r_a1</id>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p016r_1</id>", 2)
